$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.033.51"
$ws.Range("D3").Value = "3.383.56"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.96"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.97"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.66"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.83%  "
$ws.Range("E10").Value = "  -0.87%  "
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").Value = "3.962.34"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("E13").Value = "  +1.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.88"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").Value = "3.392.20"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("D17").Value = "61.121.14"
$ws.Range("E17").Value = "  +1.18%  "
$ws.Range("E18").Value = "  -1.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.60"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.37%  "
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.61"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("E22").Value = "  +3.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.550"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.22%  "
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "3.518.69"
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.188"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +6.15%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.25"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.45%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.94"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.13"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("B32").Value = "USDe"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.36"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.76%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.27"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.29%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.94"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "166.08"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").Value = "3.419.43"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.97"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.23%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.46"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.37%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0766"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.44%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.32"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.99%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.779"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.36"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.40%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.65"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.59%  "
$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.12"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.450.59"
$ws.Range("E47").Value = "  -2.47%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.92"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.63"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.10%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.12"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +9.54%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0261"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.58%  "
